$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overworld")
$ws.Range("A12").Value = ""
